$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on D-column cells whose new values look numeric,
# so Excel keeps them as text (matching the source data's inline-string type)
$ws.Range("D4:D11").NumberFormat = "@"
$ws.Range("D13:D15").NumberFormat = "@"
$ws.Range("D17:D23").NumberFormat = "@"
$ws.Range("D25:D29").NumberFormat = "@"
$ws.Range("D31:D46").NumberFormat = "@"
$ws.Range("D48:D51").NumberFormat = "@"

# Update Price (D) and Volume(1h) (E) for rows 2-31 (values refreshed by the crawler)
$ws.Range("D2").Value = '24.053.48'
$ws.Range("E2").Value = '  +16.79%  '
$ws.Range("D3").Value = '1.666.66'
$ws.Range("E3").Value = '  +12.81%  '
$ws.Range("D4").Value = '0.9989'
$ws.Range("E4").Value = '  -1.01%  '
$ws.Range("D5").Value = '309.28'
$ws.Range("E5").Value = '  +11.36%  '
$ws.Range("D6").Value = '0.9951'
$ws.Range("E6").Value = '  +3.91%  '
$ws.Range("D7").Value = '0.3723'
$ws.Range("E7").Value = '  +4.92%  '
$ws.Range("D8").Value = '0.3447'
$ws.Range("E8").Value = '  +12.03%  '
$ws.Range("D9").Value = '47.60'
$ws.Range("E9").Value = '  +20.71%  '
$ws.Range("D10").Value = '1.172'
$ws.Range("E10").Value = '  +7.97%  '
$ws.Range("D11").Value = '0.07261'
$ws.Range("E11").Value = '  +9.14%  '
$ws.Range("E12").Value = '  -0.79%  '
$ws.Range("D13").Value = '20.64'
$ws.Range("E13").Value = '  +13.69%  '
$ws.Range("D14").Value = '6.040'
$ws.Range("E14").Value = '  +10.22%  '
$ws.Range("D15").Value = '6.750'
$ws.Range("E15").Value = '  +9.24%  '
$ws.Range("D16").Value = '1.666.49'
$ws.Range("E16").Value = '  +12.94%  '
$ws.Range("D17").Value = '0.00001096'
$ws.Range("E17").Value = '  +7.77%  '
$ws.Range("D18").Value = '0.9939'
$ws.Range("E18").Value = '  +3.72%  '
$ws.Range("D19").Value = '0.06718'
$ws.Range("E19").Value = '  +12.26%  '
$ws.Range("D20").Value = '81.66'
$ws.Range("E20").Value = '  +18.02%  '
$ws.Range("D21").Value = '16.45'
$ws.Range("E21").Value = '  +13.29%  '
$ws.Range("D22").Value = '6.126'
$ws.Range("E22").Value = '  +11.33%  '
$ws.Range("D23").Value = '12.01'
$ws.Range("E23").Value = '  +7.61%  '
$ws.Range("D24").Value = '23.960.68'
$ws.Range("E24").Value = '  +16.32%  '
$ws.Range("D25").Value = '2.377'
$ws.Range("E25").Value = '  +4.09%  '
$ws.Range("D26").Value = '2.696'
$ws.Range("E26").Value = '  +28.89%  '
$ws.Range("D27").Value = '3.383'
$ws.Range("E27").Value = '  -8.10%  '
$ws.Range("D28").Value = '152.04'
$ws.Range("E28").Value = '  +2.84%  '
$ws.Range("D29").Value = '19.56'
$ws.Range("E29").Value = '  +13.92%  '
$ws.Range("D30").Value = '1.844.04'
$ws.Range("E30").Value = '  +12.77%  '
$ws.Range("D31").Value = '127.03'
$ws.Range("E31").Value = '  +10.89%  '

# Rows 32-33 and 37-40 were re-ranked (reordered) by the crawler; update Coin, Link, Price, Volume
$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").Value = '6.316'
$ws.Range("E32").Value = '  +27.43%  '
$ws.Range("B33").Value = 'HuobiToken'
$ws.Range("C33").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D33").Value = '4.120'
$ws.Range("E33").Value = '  +4.49%  '

# Rows 34-36: values refreshed
$ws.Range("D34").Value = '0.9818'
$ws.Range("E34").Value = '  +22.56%  '
$ws.Range("D35").Value = '1.734'
$ws.Range("E35").Value = '  +20.77%  '
$ws.Range("D36").Value = '0.08384'
$ws.Range("E36").Value = '  +5.87%  '

# Rows 37-40 reordered
$ws.Range("B37").Value = 'FraxShare'
$ws.Range("C37").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D37").Value = '9.039'
$ws.Range("E37").Value = '  +21.70%  '
$ws.Range("B38").Value = 'Aptos'
$ws.Range("C38").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D38").Value = '12.34'
$ws.Range("E38").Value = '  +20.05%  '
$ws.Range("B39").Value = 'InternetComputer(DFINITY)'
$ws.Range("C39").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D39").Value = '5.320'
$ws.Range("E39").Value = '  +12.75%  '
$ws.Range("B40").Value = 'Hedera'
$ws.Range("C40").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D40").Value = '0.06356'
$ws.Range("E40").Value = '  +11.61%  '

# Rows 41-51: values refreshed
$ws.Range("D41").Value = '1.292'
$ws.Range("E41").Value = '  +7.87%  '
$ws.Range("D42").Value = '0.02320'
$ws.Range("E42").Value = '  +14.93%  '
$ws.Range("D43").Value = '0.2077'
$ws.Range("E43").Value = '  +12.13%  '
$ws.Range("D44").Value = '0.6109'
$ws.Range("E44").Value = '  +16.78%  '
$ws.Range("D45").Value = '0.9941'
$ws.Range("E45").Value = '  +3.68%  '
$ws.Range("D46").Value = '3.828'
$ws.Range("E46").Value = '  +8.87%  '
$ws.Range("E47").Value = '  +9.91%  '
$ws.Range("D48").Value = '0.5959'
$ws.Range("E48").Value = '  +15.13%  '
$ws.Range("D49").Value = '127.49'
$ws.Range("E49").Value = '  +6.12%  '
$ws.Range("D50").Value = '1.997'
$ws.Range("E50").Value = '  +10.14%  '
$ws.Range("D51").Value = '0.07080'
$ws.Range("E51").Value = '  +10.53%  '
